$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Jamal Murray -> Bradley Beal (and his position/team)
$ws.Range("A2").Value = "Bradley Beal"
$ws.Range("B2").Value = "PG,SG,SF"
$ws.Range("C2").Value = "Phoenix Suns"

# Row 10: Myles Turner -> Deandre Ayton (shift up)
$ws.Range("A10").Value = "Deandre Ayton"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Portland Trail Blazers"

# Row 11: Deandre Ayton -> Alperen Sengün
$ws.Range("A11").Value = "Alperen Sengün"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Houston Rockets"

# Row 14: Obi Toppin -> Domantas Sabonis
$ws.Range("A14").Value = "Domantas Sabonis"
$ws.Range("B14").Value = "C"
$ws.Range("C14").Value = "Sacramento Kings"

# Row 15: Domantas Sabonis -> Amen Thompson
$ws.Range("A15").Value = "Amen Thompson"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Houston Rockets"

# Row 16: Bradley Beal -> Trayce Jackson-Davis
$ws.Range("A16").Value = "Trayce Jackson-Davis"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Golden State Warriors"
